# Add a new "total_staff" worksheet after the existing "annual_budget"
# sheet, populate it with header + data rows, and make it the active
# (selected/visible) sheet, matching the commit's xlsx diff.

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook and rename it.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "total_staff"

# Header row.
$newSheet.Range("A1").Value = "year"
$newSheet.Range("B1").Value = "system"
$newSheet.Range("C1").Value = "value"
$newSheet.Range("D1").Value = "staff_type"

# Data rows. Write the "SUPPORT" string before "SUPERVISION" so that the
# shared-strings table is built in the same order as the target workbook
# (SUPPORT -> index 7, SUPERVISION -> index 8).
$newSheet.Range("D3").Value = "SUPPORT"
$newSheet.Range("D2").Value = "SUPERVISION"

$newSheet.Range("A2").Value = 2021
$newSheet.Range("B2").Value = "both"
$newSheet.Range("C2").Value = 100

$newSheet.Range("A3").Value = 2021
$newSheet.Range("B3").Value = "both"
$newSheet.Range("C3").Value = 50

# Select D2 on the new sheet (matches recorded selection in the diff).
$newSheet.Range("D2").Select()

# Make the new sheet the active/selected tab.
$newSheet.Activate()
